# Refresh the NATMI ligand-receptor (Tgfb2 -> Eng) TPM output table.
# The full 5x5 sending/target-cluster grid (ECs, FAPs, Inflammatory-Mac, MuSCs,
# Resolving-Mac) now includes "MuSCs" as a sending cluster (rows 17-21 are new),
# and every metric column (E:T) is recomputed against the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Tgfb2", "Eng", "ECs", 3.0, 1.0, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3.0, 1.0, 266.9240163333333, 800.7720489999999, 0.7873936103073201, 0.78739361030732, 738.7052751780753, 6648.347476602677, 0.09389546942075302, 0.09389546942075302),
    @("ECs", "Tgfb2", "Eng", "FAPs", 3.0, 1.0, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3.0, 1.0, 51.70670833333333, 155.120125, 0.1525285446808506, 0.1525285446808506, 143.0969709380833, 1287.87273844275, 0.01818879289264614, 0.01818879289264614),
    @("ECs", "Tgfb2", "Eng", "Inflammatory-Mac", 3.0, 1.0, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3.0, 1.0, 7.729888333333332, 23.189665, 0.02280223700236483, 0.02280223700236483, 21.39226498540333, 192.53038486863, 0.002719131472688311, 0.002719131472688312),
    @("ECs", "Tgfb2", "Eng", "MuSCs", 3.0, 1.0, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3.0, 1.0, 3.652135, 10.956405, 0.01077335716164485, 0.01077335716164485, 10.10718865699, 90.96469791291, 0.001284706168158082, 0.001284706168158082),
    @("ECs", "Tgfb2", "Eng", "Resolving-Mac", 3.0, 1.0, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3.0, 1.0, 8.984181666666666, 26.952545, 0.02650225084781963, 0.02650225084781963, 24.86348917377667, 223.77140256399, 0.0031603523974386, 0.003160352397438601),
    @("FAPs", "Tgfb2", "Eng", "ECs", 3.0, 1.0, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3.0, 1.0, 266.9240163333333, 800.7720489999999, 0.7873936103073201, 0.78739361030732, 3439.523795289052, 30955.71415760146, 0.43719154606641, 0.4371915460664099),
    @("FAPs", "Tgfb2", "Eng", "FAPs", 3.0, 1.0, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3.0, 1.0, 51.70670833333333, 155.120125, 0.1525285446808506, 0.1525285446808506, 666.2811991652222, 5996.530792487, 0.08468977827017633, 0.08468977827017633),
    @("FAPs", "Tgfb2", "Eng", "Inflammatory-Mac", 3.0, 1.0, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3.0, 1.0, 7.729888333333332, 23.189665, 0.02280223700236483, 0.02280223700236483, 99.60563018138221, 896.4506716324399, 0.0126606885277437, 0.0126606885277437),
    @("FAPs", "Tgfb2", "Eng", "MuSCs", 3.0, 1.0, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3.0, 1.0, 3.652135, 10.956405, 0.01077335716164485, 0.01077335716164485, 47.06060327078666, 423.54542943708, 0.005981786760990885, 0.005981786760990885),
    @("FAPs", "Tgfb2", "Eng", "Resolving-Mac", 3.0, 1.0, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3.0, 1.0, 8.984181666666666, 26.952545, 0.02650225084781963, 0.02650225084781963, 115.7681764577911, 1041.91358812012, 0.01471508007015176, 0.01471508007015176),
    @("Inflammatory-Mac", "Tgfb2", "Eng", "ECs", 1.0, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3.0, 1.0, 266.9240163333333, 800.7720489999999, 0.7873936103073201, 0.78739361030732, 9.137253952450555, 82.23528557205499, 0.001161419551085723, 0.001161419551085723),
    @("Inflammatory-Mac", "Tgfb2", "Eng", "FAPs", 1.0, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3.0, 1.0, 51.70670833333333, 155.120125, 0.1525285446808506, 0.1525285446808506, 1.770006804097222, 15.930061236875, 0.0002249823107173179, 0.0002249823107173179),
    @("Inflammatory-Mac", "Tgfb2", "Eng", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3.0, 1.0, 7.729888333333332, 23.189665, 0.02280223700236483, 0.02280223700236483, 0.2646069607972222, 2.381462647175, 0.00003363370430793883, 0.00003363370430793884),
    @("Inflammatory-Mac", "Tgfb2", "Eng", "MuSCs", 1.0, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3.0, 1.0, 3.652135, 10.956405, 0.01077335716164485, 0.01077335716164485, 0.1250186679416667, 1.125168011475, 0.0000158908930356701, 0.0000158908930356701),
    @("Inflammatory-Mac", "Tgfb2", "Eng", "Resolving-Mac", 1.0, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3.0, 1.0, 8.984181666666666, 26.952545, 0.02650225084781963, 0.02650225084781963, 0.3075435120861111, 2.767891608775, 0.00003909129040356621, 0.00003909129040356622),
    @("MuSCs", "Tgfb2", "Eng", "ECs", 3.0, 1.0, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3.0, 1.0, 266.9240163333333, 800.7720489999999, 0.7873936103073201, 0.78739361030732, 2007.307573733051, 18065.76816359746, 0.2551451752690714, 0.2551451752690714),
    @("MuSCs", "Tgfb2", "Eng", "FAPs", 3.0, 1.0, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3.0, 1.0, 51.70670833333333, 155.120125, 0.1525285446808506, 0.1525285446808506, 388.8419958460833, 3499.57796261475, 0.04942499120731082, 0.04942499120731083),
    @("MuSCs", "Tgfb2", "Eng", "Inflammatory-Mac", 3.0, 1.0, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3.0, 1.0, 7.729888333333332, 23.189665, 0.02280223700236483, 0.02280223700236483, 58.12988883036332, 523.1689994732699, 0.007388783297624878, 0.007388783297624879),
    @("MuSCs", "Tgfb2", "Eng", "MuSCs", 3.0, 1.0, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3.0, 1.0, 3.652135, 10.956405, 0.01077335716164485, 0.01077335716164485, 27.46458841171, 247.18129570539, 0.003490973339460217, 0.003490973339460217),
    @("MuSCs", "Tgfb2", "Eng", "Resolving-Mac", 3.0, 1.0, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3.0, 1.0, 8.984181666666666, 26.952545, 0.02650225084781963, 0.02650225084781963, 67.56235782385666, 608.06122041471, 0.008587727089825702, 0.008587727089825704)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c = $c + 1
    }
    $r = $r + 1
}
